$d = $word.ActiveDocument

# Locate the block of letter-body paragraphs that need the |e (escape)
# filter added to every docassemble interpolation -- from the opening
# "{{ today() }}" date line through the closing "{{ user }}"
# signature line.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($startPara -eq $null -and $t.Contains("{{ today")) {
        $startPara = $i
    }
    if ($t.Contains("{{ user }}")) {
        $endPara = $i
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the letter-body paragraph block to update."
}

$start = $d.Paragraphs($startPara).Range.Start
$end = $d.Paragraphs($endPara).Range.End
$r = $d.Range($start, $end)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="5040"/></w:pPr><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ today</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>()</w:t></w:r><w:r><w:t>|e</w:t></w:r><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p><w:p><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ adverse</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_party.address_block()</w:t></w:r><w:r><w:t>|e</w:t></w:r><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p><w:p><w:r><w:t>To whom it may concern:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>This letter relates to the critical issue of {</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ subject</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_matter</w:t></w:r><w:r><w:t>|e</w:t></w:r><w:r><w:t xml:space="preserve"> }}.  Please cease and desist your officious intermeddling.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>If you have any questions, you can reach me at {</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ phone</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_number</w:t></w:r><w:r><w:t>|e</w:t></w:r><w:r><w:t xml:space="preserve"> }}.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="5040"/></w:pPr><w:r><w:t>Sincerely,</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="5040"/></w:pPr><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{ user.signature</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="5040"/></w:pPr><w:r><w:t>{{ user</w:t></w:r><w:r><w:t>|e</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> }}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Host "Paragraphs replaced:" $startPara "to" $endPara
